$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from "Checklist" to "Pediatrics"
$ws.Name = "Pediatrics"

# Remove the oldest log entry (old row 2 - Student 201317 / Selection at 18:04:18).
# This shifts every subsequent row up by one, preserving their original
# cell formatting/types (dates, etc.) untouched.
$ws.Rows.Item(2).Delete() | Out-Null

# The rows that used to be "Selection" scans (now rows 2-7) were re-logged
# as "Scan" entries.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 5).Value = "Scan"
}

# Row 8 (former row 9, the Manual entry) is unchanged and already reads
# "Manual" after the shift, so nothing further to do there.
